$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 3 "Ethical Considerations": append a period to the first bullet.
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = "NBA data is publicly available."

# ---------------------------------------------------------------------------
# 2. Slide 4 "React, Recharts, & MUI": duplicate it first (before resizing)
#    so the new "Data Cleaning & Storage" slide inherits the original,
#    un-resized placeholder geometry. The duplicate is inserted right after
#    slide 4, which is exactly where the new slide belongs.
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$newSlide = $s4.Duplicate()

# Append a period to the "Recharts..." bullet.
$s4.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2,1).Runs(1,1).Text = "Recharts is a way to use Reacts library to build out graphs."

# Reposition the title and content placeholders (values are EMU/12700,
# nudged by a few ULPs so the lossy single-precision round-trip inside the
# COM layer still lands on the exact target EMU).
$title4 = $s4.Shapes.Item(1)
$title4.Left = 70.4596062992126
$title4.Top = 225.05913545826772
$title4.Width = 353.28001409999996
$title4.Height = 89.88165674330708

$content4 = $s4.Shapes.Item(2)
$content4.Left = 536.2604724409449
$content4.Top = 145.54047394094488
$content4.Width = 379.2
$content4.Height = 248.91905981811024

# ---------------------------------------------------------------------------
# 3. New slide (now at index 5, id 264): set title & body text for
#    "Data Cleaning & Storage".
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "Data Cleaning & Storage"

$body5 = $s5.Shapes.Item(2).TextFrame.TextRange
$body5.Paragraphs(1,1).Runs(1,1).Text = "The CSVs were imported into a Python file for cleaning."
$body5.Paragraphs(2,1).Runs(1,1).Text = "Dropped unnecessary columns, formatted columns, etc."
$body5.Paragraphs(3,1).Runs(1,1).Text = "Data was exported into JSONs and added into MongoDB to be stored."
